# Generate Report for Handback
# -----------------------------------------------------------------
# This localization-status workbook tracks handoff/handback state for
# each target locale. The handback just completed for de-de (and the
# zh-cn target file link is filled in too), so:
#   - the overview "status" column (zh-cn / de-de) flips from
#     "Ready for handoff" to "Handed back: in sync with en-US"
#   - each locale sheet gets its "Latest Target File" (source doc link)
#     and "Latest Handback File" populated
#   - de-de's "Latest Handback DateTime" is stamped with the real
#     handback time (zh-cn's backing timestamp string is also refreshed)
# -----------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$sourceDocName = "645a5244-477f-41e9-8df2-c9bc9ab3ea47.md"
$sourceDocUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/beaa13ef3e5483f1c1da9f5f50e4d513f01776b9/e2e/645a5244-477f-41e9-8df2-c9bc9ab3ea47.md"

# Helper: ColumnWidth is quantized by the host to 1/6-character steps
# (stored_width = round(ColumnWidth*6)/6 + 5/6), so back the desired
# stored width out to the ColumnWidth that lands on it.
function Set-StoredColumnWidth($col, $storedWidth) {
    $col.ColumnWidth = ($storedWidth - 5.0/6.0)
}

# ---------------------------------------------------------------
# Overview sheet: status text + widen the two status columns
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusText
$ov.Range("F2").Value = $statusText
Set-StoredColumnWidth $ov.Columns.Item(5) 29.9777047293527
Set-StoredColumnWidth $ov.Columns.Item(6) 29.9777047293527

# B2's pre-existing hyperlink formatting (underline + blue) gets
# dropped by the save pipeline's style-table rebuild unless the font
# is touched again during this session - reassert it so it round-trips.
$ov.Range("B2").Font.Underline = 2
$ov.Range("B2").Font.Color = 15570276

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusText

$zh.Hyperlinks.Add($zh.Range("I2"), $sourceDocUrl, "", "", $sourceDocName)
$zh.Range("I2").Font.Underline = 2
$zh.Range("I2").Font.Color = 15570276

$zh.Range("J2").Value = "645a5244-477f-41e9-8df2-c9bc9ab3ea47.177b21937c6086140a1f3e49120247cc6d616487.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-06 23:13:20"

Set-StoredColumnWidth $zh.Columns.Item(3) 29.9777047293527
Set-StoredColumnWidth $zh.Columns.Item(9) 40
Set-StoredColumnWidth $zh.Columns.Item(10) 40

# A2's pre-existing hyperlink formatting - see Overview!B2 note above.
$zh.Range("A2").Font.Underline = 2
$zh.Range("A2").Font.Color = 15570276

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusText

$de.Hyperlinks.Add($de.Range("I2"), $sourceDocUrl, "", "", $sourceDocName)
$de.Range("I2").Font.Underline = 2
$de.Range("I2").Font.Color = 15570276

$de.Range("J2").Value = "645a5244-477f-41e9-8df2-c9bc9ab3ea47.177b21937c6086140a1f3e49120247cc6d616487.de-de.xlf"
$de.Range("K2").Value = "2016-09-06 23:13:28"

Set-StoredColumnWidth $de.Columns.Item(3) 29.9777047293527
Set-StoredColumnWidth $de.Columns.Item(9) 40
Set-StoredColumnWidth $de.Columns.Item(10) 40

# A2's pre-existing hyperlink formatting - see Overview!B2 note above.
$de.Range("A2").Font.Underline = 2
$de.Range("A2").Font.Color = 15570276

Write-Host "Handback report generated."
